$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (QA) - update existing claim number
$ws.Range("A2").Value = "QA"
$ws.Range("B2").Value = "'0420194406717"
$ws.Range("C2").Value = 100

# Row 3 (QA) - new row
$ws.Range("A3").Value = "QA"
$ws.Range("B3").Value = "'1120194100412"
$ws.Range("C3").Value = 100

# Row 4 (QA) - new row
$ws.Range("A4").Value = "QA"
$ws.Range("B4").Value = "'1220194200667"
$ws.Range("C4").Value = 100

# Row 5 (PREPROD) - moved from old row 3
$ws.Range("A5").Value = "PREPROD"
$ws.Range("B5").Value = "'1120170200928 "
$ws.Range("C5").Value = 100

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("B5").Select() | Out-Null
